# Update odds values on Sheet1 to match the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 1.95
$ws.Range("S2").Value = 1.57

$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38

$ws.Range("G6").Value = 2.35
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 2.75
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 3.3
$ws.Range("V6").Value = 2.07
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 25
$ws.Range("AD6").Value = 6.6
$ws.Range("AE6").Value = 12.5
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 15
$ws.Range("AT6").Value = 2.62
$ws.Range("AU6").Value = 6.7
$ws.Range("AW6").Value = 4.75

$ws.Range("S7").Value = 1.27
$ws.Range("T7").Value = 3.54

$ws.Range("M13").Value = 1.05
$ws.Range("O13").Value = 1.3

$ws.Range("M14").Value = 1.05
$ws.Range("O14").Value = 1.3

$ws.Range("M15").Value = 1.07
$ws.Range("O15").Value = 1.41
$ws.Range("P15").Value = 2.62

$ws.Range("M16").Value = 1.05
$ws.Range("O16").Value = 1.3

$ws.Range("Q17").Value = 1.65
$ws.Range("R17").Value = 2.2
